# Deepthi - First Test case - View Licence Details
#
# Add a new "Galaxy S7 Edge " worksheet, placed right after the
# existing "Galaxy S7 " sheet, with the same "Test Parameters" / value
# header row used by the other device sheets.
#
# We create it by copying the "Galaxy S7 " sheet (which already has the
# exact A1/B1 layout, column widths and shared-string-backed cell values
# we need) rather than building it cell-by-cell, so the new sheet ends up
# with the same cell types (shared string text, not numbers) and styling
# as its siblings.

$wb = $excel.ActiveWorkbook

$sourceSheet = $wb.Worksheets.Item("Galaxy S7 ")
$sourceSheet.Copy([System.Reflection.Missing]::Value, $sourceSheet)

$newSheet = $wb.ActiveSheet
$newSheet.Name = "Galaxy S7 Edge "

# Restore the originally active sheet/tab selection so the only
# meaningful change is the addition of the new worksheet.
$wb.Worksheets.Item(1).Activate()
